$d = $word.ActiveDocument

$old = "Stanford GSB (scheduled), Princeton (scheduled), Kellogg (scheduled), Stanford GSB Junior Faculty Workshop on Financial Regulation and Banking (scheduled), London Business School Female Economist Conference (scheduled), IMF (scheduled), Wisconsin Money, Banking, and Asset Markets Conference (scheduled)"

$new = "Stanford GSS, Princeton, Kellogg, Stanford GSB Junior Faculty Workshop on Financial Regulation and Banking (scheduled), London Business School Female Economist Conference, IMF, Wisconsin Money, Banking, and Asset Markets Conference, Finance Theory Group (Boston College), SFS (Toronto), FIRS (Lisbon), WFA (Park City)"

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the target conference-presentations text to replace."
}
